$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# The legend grows from 5 rows to 7: two new status rows are inserted between
# the existing "Class ID / DLL missing" row (row 3) and the existing
# "No Class ID" row (old row 4). Rather than retyping the two rows that only
# shift position (which contain a curly apostrophe / en dash that must stay
# byte-identical), copy their value+format down to rows 6/7 first, then
# overwrite rows 4/5 with the two brand-new legend entries.
# ---------------------------------------------------------------------------

# Old row 5 ("Sparx key doesn't exist...", red) -> row 7
$ws.Range("A5").Copy()
$ws.Range("A7").PasteSpecial(-4163)   # xlPasteValues
$ws.Range("A5").Copy()
$ws.Range("A7").PasteSpecial(-4122)   # xlPasteFormats

# Old row 4 ("...No Class ID is set...", hot pink) -> row 6
$ws.Range("A4").Copy()
$ws.Range("A6").PasteSpecial(-4163)
$ws.Range("A4").Copy()
$ws.Range("A6").PasteSpecial(-4122)

$ws.Range("A6").RowHeight = 15.75
$ws.Range("A7").RowHeight = 15.75

# New row 4: "Mismatch on AddIn name" (yellow fill). Base the border/font on
# the already-boxed row 3 style, then restyle just the fill + text.
$ws.Range("A3").Copy()
$ws.Range("A4").PasteSpecial(-4122)
$ws.Range("A4").Value = "Sparx key exist and Class ID - Mismatch on AddIn name "
$ws.Range("A4").Interior.Color = 65535

# New row 5: "DLL does not appear to be a normal AddIn" (pale green fill)
$ws.Range("A3").Copy()
$ws.Range("A5").PasteSpecial(-4122)
$ws.Range("A5").Value = "Sparx key exist and Class ID - DLL does not appear to be a normal AddIn (integration?)"
$ws.Range("A5").Interior.Color = 10025880

# ---------------------------------------------------------------------------
# Column C: plain-text hex reference of the fill color used by each legend
# row (kept in the same row order as column A).
# ---------------------------------------------------------------------------
$ws.Range("C1").Value = "00FF00"
$ws.Range("C2").Value = "F5DE83"
$ws.Range("C3").Value = "FFC0CB"
$ws.Range("C4").Value = "FFFF00"
$ws.Range("C5").Value = "98FB98"
$ws.Range("C6").Value = "FF69B4"
$ws.Range("C7").Value = "FF0000"

# Refresh the selection to the full, newly-sized legend column.
$ws.Range("A1:A7").Select()
